# Generate Report for Handback
#
# Two of the tracked files ("9eaafee2-...md" and "cc0121d2-...md") have now
# been handed back and are in sync with en-US. This updates:
#   - the Overview sheet's Status column for both files
#   - each locale sheet's Status column, newly-populated "Latest Target
#     File" / "Latest Handback File" hyperlinks, and "Latest Handback
#     DateTime" for both files

$wb = $excel.ActiveWorkbook

$newStatus = "Handed back: in sync with en-US"

# Cornflower blue (FF6495ED) packed as an OLE BGR integer (R + G*256 + B*65536),
# matching the existing hyperlink-style cells (A2/A3/C2/C3, ...).
$hyperlinkColor = 15570276

function Set-HandoffHyperlink($ws, $cellRef, $address, $displayText) {
    $ws.Hyperlinks.Add($ws.Range($cellRef), $address, "", "", $displayText) | Out-Null
    $font = $ws.Range($cellRef).Font
    $font.Color = $hyperlinkColor
    $font.Underline = 2
}

# ---------------------------------------------------------------------
# Overview sheet: both tracked files show the new status in zh-cn/de-de
# ---------------------------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("B2").Value = $newStatus
$overview.Range("C2").Value = $newStatus
$overview.Range("B3").Value = $newStatus
$overview.Range("C3").Value = $newStatus

# ---------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------
$zhcn = $wb.Worksheets.Item("zh-cn")

$zhcn.Range("B2").Value = $newStatus
$zhcn.Range("B3").Value = $newStatus

Set-HandoffHyperlink $zhcn "E2" `
    "https://github.com/OpenLocalizationTest/oltest/blob/cfe8cbc8df5900471171ea8c11453f19c0ca7f2c/e2e/9eaafee2-ed77-4467-a92d-a45e487e0bb0.md" `
    "9eaafee2-ed77-4467-a92d-a45e487e0bb0.md"

Set-HandoffHyperlink $zhcn "F2" `
    "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/8eb2d5370fe60223dfe4bb2add4261effefea5ba/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/9eaafee2-ed77-4467-a92d-a45e487e0bb0.ff53143d17f478be72e13c2a353edef2e110344a.zh-cn.xlf" `
    "9eaafee2-ed77-4467-a92d-a45e487e0bb0.ff53143d17f478be72e13c2a353edef2e110344a.zh-cn.xlf"

$zhcn.Range("G2").Value = "2016-03-10 07:08:50"

Set-HandoffHyperlink $zhcn "E3" `
    "https://github.com/OpenLocalizationTest/oltest/blob/cfe8cbc8df5900471171ea8c11453f19c0ca7f2c/e2e/cc0121d2-f0d0-4fec-8829-8fd6ee143b4b.md" `
    "cc0121d2-f0d0-4fec-8829-8fd6ee143b4b.md"

Set-HandoffHyperlink $zhcn "F3" `
    "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/8eb2d5370fe60223dfe4bb2add4261effefea5ba/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/cc0121d2-f0d0-4fec-8829-8fd6ee143b4b.68c5c99a7a296575aafb481b4ea1c24729e71c72.zh-cn.xlf" `
    "cc0121d2-f0d0-4fec-8829-8fd6ee143b4b.68c5c99a7a296575aafb481b4ea1c24729e71c72.zh-cn.xlf"

$zhcn.Range("G3").Value = "2016-03-10 07:08:50"

# ---------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------
$dede = $wb.Worksheets.Item("de-de")

$dede.Range("B2").Value = $newStatus
$dede.Range("B3").Value = $newStatus

Set-HandoffHyperlink $dede "E2" `
    "https://github.com/OpenLocalizationTest/oltest/blob/cfe8cbc8df5900471171ea8c11453f19c0ca7f2c/e2e/9eaafee2-ed77-4467-a92d-a45e487e0bb0.md" `
    "9eaafee2-ed77-4467-a92d-a45e487e0bb0.md"

Set-HandoffHyperlink $dede "F2" `
    "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/025fdafac48e9ccf229d5444711d7735798218cd/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/9eaafee2-ed77-4467-a92d-a45e487e0bb0.ff53143d17f478be72e13c2a353edef2e110344a.de-de.xlf" `
    "9eaafee2-ed77-4467-a92d-a45e487e0bb0.ff53143d17f478be72e13c2a353edef2e110344a.de-de.xlf"

$dede.Range("G2").Value = "2016-03-10 07:09:02"

Set-HandoffHyperlink $dede "E3" `
    "https://github.com/OpenLocalizationTest/oltest/blob/cfe8cbc8df5900471171ea8c11453f19c0ca7f2c/e2e/cc0121d2-f0d0-4fec-8829-8fd6ee143b4b.md" `
    "cc0121d2-f0d0-4fec-8829-8fd6ee143b4b.md"

Set-HandoffHyperlink $dede "F3" `
    "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/025fdafac48e9ccf229d5444711d7735798218cd/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/cc0121d2-f0d0-4fec-8829-8fd6ee143b4b.68c5c99a7a296575aafb481b4ea1c24729e71c72.de-de.xlf" `
    "cc0121d2-f0d0-4fec-8829-8fd6ee143b4b.68c5c99a7a296575aafb481b4ea1c24729e71c72.de-de.xlf"

$dede.Range("G3").Value = "2016-03-10 07:09:02"
